# Add five new "deepvein" (vent) resource nodes to the
# rebalance_localizations sheet, keeping the A2:K... table sorted
# alphabetically by column A (resource_name key), exactly like the
# surrounding rows are already sorted.
#
# New rows (final sheet position -> resource key / English label):
#   432 -> resource_name/fluorine_deepvein       / Fluorine Vent
#   440 -> resource_name/nitric_acid_deepvein    / Nitric acid well
#   446 -> resource_name/plasma_charged_deepvein / Charged Plasma Vent
#   447 -> resource_name/plasma_deepvein         / Plasma Vent
#   449 -> resource_name/resin_deepvein          / Resin well

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rebalance_localizations")

# Insert the 5 blank rows first (bottom-to-top, using the ORIGINAL
# row numbers so each Insert() target is still valid for the rows
# that have not been processed yet).
$ws.Rows.Item(445).Insert()   # will hold resin_deepvein   (-> row 449)
$ws.Rows.Item(444).Insert()   # will hold plasma_deepvein  (-> row 447)
$ws.Rows.Item(444).Insert()   # will hold plasma_charged   (-> row 446)
$ws.Rows.Item(439).Insert()   # will hold nitric_acid_dv   (-> row 440)
$ws.Rows.Item(432).Insert()   # will hold fluorine_dv      (-> row 432)

# Column A (resource_name keys)
$ws.Range("A446").Value = "resource_name/plasma_charged_deepvein"
$ws.Range("A447").Value = "resource_name/plasma_deepvein"
$ws.Range("A432").Value = "resource_name/fluorine_deepvein"
$ws.Range("A440").Value = "resource_name/nitric_acid_deepvein"
$ws.Range("A449").Value = "resource_name/resin_deepvein"

# Column B (English display names)
$ws.Range("B432").Value = "Fluorine Vent"
$ws.Range("B446").Value = "Charged Plasma Vent"
$ws.Range("B447").Value = "Plasma Vent"
$ws.Range("B440").Value = "Nitric acid well"
$ws.Range("B449").Value = "Resin well"

# Restore the view state seen in the edited workbook: unfrozen-pane
# scroll back to the top and the active cell on B14.
$ws.Range("B14").Select() | Out-Null

Write-Output "Inserted 5 deepvein/vent resource rows; sheet now spans A1:K458"
